$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 20:22"

# Update Alemania stats (row 9)
$ws.Cells.Item(9,2).Value = 160943
$ws.Cells.Item(9,3).Value = 1031
$ws.Cells.Item(9,5).Value = 34167
$ws.Cells.Item(9,7).Value = 62
$ws.Cells.Item(9,8).Value = 6376

# Update Emiratos Arabes Unidos stats (row 37)
$ws.Cells.Item(37,2).Value = 11929
$ws.Cells.Item(37,3).Value = 549
$ws.Cells.Item(37,4).Value = 2329
$ws.Cells.Item(37,5).Value = 9502
$ws.Cells.Item(37,7).Value = 9
$ws.Cells.Item(37,8).Value = 98

# Update Montenegro stats (row 125) - only Muertes hoy (F) changes
$ws.Cells.Item(125,6).Value = 2

# Update Angola stats (row 180) - Casos activos (D) and Recuperados (E) change
$ws.Cells.Item(180,4).Value = 7
$ws.Cells.Item(180,5).Value = 18

# Suazilandia's case counts increased enough to move it up in the
# (descending, by total cases) sort order: it used to sit between
# Guinea-Bisau and Benin (row 163); now it belongs between Monaco and
# Liechtenstein (row 153). Remove the old row and insert a new one in
# its new sorted position with the updated figures.
$ws.Rows("163:163").Delete()
$ws.Rows("153:153").Insert()
$ws.Cells.Item(153,1).Value = "Suazilandia"
$ws.Cells.Item(153,2).Value = 91
$ws.Cells.Item(153,3).Value = 20
$ws.Cells.Item(153,4).Value = 10
$ws.Cells.Item(153,5).Value = 80
$ws.Cells.Item(153,6).Value = 0
$ws.Cells.Item(153,7).Value = 0
$ws.Cells.Item(153,8).Value = 1
